# Apply crypto price/volume updates per the commit diff.
# All data cells in this sheet are stored as text (prices use
# "." as a thousands separator and volumes carry surrounding
# spaces/sign), so every write keeps the cell a string. Values
# that look like plain numbers (e.g. "1.002") get a leading
# apostrophe -- the standard Excel "force text" marker -- so
# they are not auto-converted to numeric cells; values that
# already fail to parse as a number (e.g. "27.663.83", a URL,
# or a coin name) are written as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.663.83"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.878.91"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'331.34"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4717"
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").Value = "'48.25"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").Value = "'0.08031"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "'1.023"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "'21.76"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "1.874.26"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "'5.963"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'7.166"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "'87.05"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'0.06621"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "'17.18"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "27.675.31"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'5.507"
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'2.293"
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("D26").Value = "2.093.28"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'155.86"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "'20.28"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("D29").Value = "'2.090"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'5.590"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").Value = "'122.27"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.9672"
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.09564"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("D35").Value = "'3.632"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'5.304"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "'0.06107"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "'0.02254"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").Value = "'1.228"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").Value = "'8.107"
$ws.Range("E40").Value = "  -4.80%  "
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "'0.1902"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "'1.254"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "'0.5696"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'12.19"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").Value = "'1.932"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "'0.06823"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'111.27"
$ws.Range("E51").Value = "  +1.63%  "
